$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 9, shifting old rows 9-13 down to rows 10-14.
# This grows the table from 13 to 14 data rows.
$ws.Rows.Item(9).EntireRow.Insert()

# Populate the new / changed "Matrix tutorial" instruction rows (6-9).
# Values are set in the same order the strings were introduced into the
# shared-string table by the original author so the resulting shared
# string indices line up with the target workbook.
$ws.Range("B6").Value = "במהלך הניסוי יוצגו לך פרופילים של אנשים שונים. לאחר שתתבונן בפרופיל של כל אדם, תשחק מול משחק הקשור לסיואציה בו אתם כביכול נפגשים. המשחק מתאר שתי בחירות אפשריות שלך  ושתי בחירות אפשריות של האחר. עבור כל שילוב של שתי בחירות ישנה תוצאה אשר מיוצגת במספר, תוצאה עבורך ותוצאה עבור האחר. לחץ על הכפתור Ready to play על מנת להתבונן במשחק"
$ws.Range("B8").Value = "יפה מאוד! כאשר אתה לוחץ על אחת השורות, זה אומר שבחרת באפשרות הזו. כמה שתקבל תלוי במה יבחר השחקן השני. בזמן שאתה מבצע את בחירתך אינך יכול לראות מה השחקן השני בחר, אך גם כשאר השחקן השני בוחר, אתה לא יכול לראות את הבחירה שלו. תוצאות המשחק יקבעו רק לאחר מכן. בוא נמשיך בלמידה - לחץ עכשיו על האופציה העליונה."
$ws.Range("B9").Value = "מעולה! כעת, ענה על השאלות שמופיעות מטה. בכל שאלה תצטרך לומר כמה אתה או השחק האחר יקבלו עבור שילוב מסוים של שתי אפשרויות, הבחירה שלך והבחירה של האחר."
$ws.Range("B7").Value = "במשחק זה אתה יכול לבחור בין {} , לבין {} וכך גם האדם השני. כל אחד מכם בוחר לפי ראות עיניו, ומתוך הניסיון להגיע לתוצאה הטובה ביותר עבורו. אתה יכול לבחור או בשורה הראשונה ({}) או בשורה השניה ({}). כך גם השחקן השני, אל שהוא יכול לבחור בין הטור הימני לשמאלי. בכל צירוף של שתי בחירות רשום כמה אתה תקבל וכמה האחר. התוצאה שלך רשומה תמיד מצד שמאל. לצורך בדיקת הבנה - בחר כעת את השורה התחונה - {}."

# The new row 9 needs its A / C columns filled in (copied from the other
# "Matrix tutorial" / "irrelevant" rows).
$ws.Range("A9").Value = "Matrix tutorial"
$ws.Range("C9").Value = "irrelevant"

# Update the active cell selection to match the saved workbook state.
[void]$ws.Range("B9").Select()
